$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): rename the team-member columns D:H and add column H
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Alex"
$ws.Range("E1").Value = "Georgios"
$ws.Range("F1").Value = "Karen"
$ws.Range("G1").Value = "Roshi"
$ws.Range("H1").Value = "Stuart"

# ---------------------------------------------------------------------------
# Row 2 - meeting 1
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = 43866
$ws.Range("C2").Value = "Group Formation: set up communication channel in Slack and GitHub repository"
$ws.Range("D2:H2").Value = "yes"

# ---------------------------------------------------------------------------
# Row 3 - meeting 2
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = 43872
$ws.Range("C3").Value = 'Agreed topic of "Plastic Pollution", distributed research activity for week'
$ws.Range("D3:H3").Value = "yes"

# ---------------------------------------------------------------------------
# Row 4 - meeting 3 (date is now a formula)
# ---------------------------------------------------------------------------
$ws.Range("B4").Formula = "=B3+7"
$ws.Range("C4").Value = 'Presented inividuals'' research findings and discussed hypothesis'
$ws.Range("D4:H4").Value = "yes"

# ---------------------------------------------------------------------------
# Row 5 - meeting 4
# ---------------------------------------------------------------------------
$ws.Range("B5").Formula = "=B4+7"
$ws.Range("C5").Value = 'Decided on final dataset to use and hypothesis of "proportion of marine plastics pollution does not change over time"'
$ws.Range("D5:H5").Value = "yes"

# ---------------------------------------------------------------------------
# Row 6 - meeting 5
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = 43894
$ws.Range("C6").Value = "Presentation draft agreed"
$ws.Range("D6:H6").Value = "yes"

# ---------------------------------------------------------------------------
# Row 7 - meeting 6
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = 43900
$ws.Range("C7").Value = "Distributed section writing activity for week"
$ws.Range("D7:H7").Value = "yes"

# ---------------------------------------------------------------------------
# Rows 8-9 lose their Topic/attendance cells, only No/Date remain
# ---------------------------------------------------------------------------
$ws.Range("C8:G8").ClearContents()
$ws.Range("B8").Formula = "=B7+7"

$ws.Range("C9:G9").ClearContents()
$ws.Range("B9").Formula = "=B8+7"

# ---------------------------------------------------------------------------
# New rows 10-13 (meetings 9-12), only No/Date populated
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = 9
$ws.Range("B10").Formula = "=B9+7"

$ws.Range("A11").Value = 10
$ws.Range("B11").Formula = "=B10+7"

$ws.Range("A12").Value = 11
$ws.Range("B12").Formula = "=B11+7"

$ws.Range("A13").Value = 12
$ws.Range("B13").Formula = "=B12+7"

# ---------------------------------------------------------------------------
# Rows 14-16: empty, date-formatted placeholder cells in column B
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Column widths / selection cosmetics
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 25.8
$ws.Range("C14").Select() | Out-Null
